# Apply "Added altimeter config and basic readouts" edit to the Translation sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New rows appended after the existing last row (151), mirroring the
# "Debug SD card" block pattern for a new "Debug EnvSensors" block, plus
# new altimeter readouts (Temperature / Pressure / Altitude).
$rows = @(
    @("SingleUseId199", "Default", "Center", "LTR", "Debug EnvSensors"),
    @("SingleUseId200", "Small",   "Left",   "LTR", "Output files: <value>"),
    @("SingleUseId201", "Small",   "Left",   "LTR", "0"),
    @("SingleUseId202", "Small",   "Left",   "LTR", "Input files: <value>"),
    @("SingleUseId203", "Small",   "Left",   "LTR", "0"),
    @("SingleUseId204", "Small",   "Left",   "LTR", "Free space: <value>"),
    @("SingleUseId205", "Small",   "Left",   "LTR", "0"),
    @("SingleUseId206", "Small",   "Left",   "LTR", "Total space: <value>"),
    @("SingleUseId207", "Small",   "Left",   "LTR", "0"),
    @("SingleUseId208", "Small",   "Left",   "LTR", "State: <value>"),
    @("SingleUseId209", "Small",   "Left",   "LTR", "UNINITIALIZED"),
    @("SingleUseId214", "Small",   "Left",   "LTR", "Temperature: <value> C"),
    @("SingleUseId215", "Small",   "Left",   "LTR", "0.00"),
    @("SingleUseId216", "Small",   "Left",   "LTR", "Pressure: <value> hPa"),
    @("SingleUseId217", "Small",   "Left",   "LTR", "0.00"),
    @("SingleUseId218", "Small",   "Left",   "LTR", "Altitude: <value> m"),
    @("SingleUseId219", "Small",   "Left",   "LTR", "0.00")
)

$startRow = 152
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}

# Columns F on rows 154/156/158/160 hold the purely-numeric-looking text
# "0" (matching existing text id rows such as F41/F43/F45/F47). Assigning
# ".Value = '0'" directly would be auto-coerced to a numeric constant, so
# instead copy the text value from an existing "0" text cell to preserve
# the text data type without altering cell formatting.
$zeroSource = $ws.Range("F41")
$zeroSource.Copy() | Out-Null
foreach ($r in @(154, 156, 158, 160)) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = $false

# Columns F on rows 164/166/168 hold the text "0.00" (matching existing
# text id rows such as F49/F61/F63/F72).
$decimalSource = $ws.Range("F49")
$decimalSource.Copy() | Out-Null
foreach ($r in @(164, 166, 168)) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = $false
